$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: column headers - translate slug-style identifiers into
# human-readable Spanish labels.
$ws.Range("A1").Value = "CCAA 2ª residencia, código"
$ws.Range("B1").Value = "Comarca nombre"
$ws.Range("C1").Value = "Número hogares"
$ws.Range("D1").Value = "Comarca código"
$ws.Range("E1").Value = "CCAA 2ª residencia, nombre"
$ws.Range("F1").Value = "Aragón"
$ws.Range("G1").Value = "Municipio código"
$ws.Range("H1").Value = "Municipio nombre"

# Row 2: DSD concept reference (measure/dimension URI) per column.
$ws.Range("A2").Value = "null"
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("C2").Value = "iaest-measure:numero-hogares"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "iaest-measure:ccaa-2-residencia-nombre"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = "sdmx-dimension:refArea"

# Row 3: role of the column in the DSD (medida/dim).
$ws.Range("A3").Value = "null"
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "dim"
$ws.Range("G3").Value = "null"
$ws.Range("H3").Value = "dim"

# Row 4: data type / codelist URI template for the column.
$ws.Range("A4").Value = "null"
$ws.Range("B4").Value = "URI-comarca"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "xsd:string"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "null"
$ws.Range("H4").Value = "URI-Municipio"
